# epitweetr subscribers.xlsx - add "one tweet alert" related columns
#
# 1) Adds three new trailing columns (K:M) to the subscribers sheet:
#      K - "One tweet alerts"
#      L - "Topics ignoring 1 tweet alerts "
#      M - "Regions ignoring 1 tweet alerts "
#    These reuse the same header/body formatting already used by the
#    existing columns (B:J) so no new styles are introduced beyond what
#    wrapping the header text requires.
# 2) Turns wrapping on for the header row (B2:M2) and grows its row
#    height to fit the now-taller, wrapped header captions.
# 3) Slightly widens the two brand new columns and nudges the sheet's
#    default column width.
# 4) Leaves the saved cursor/selection on D16, matching the authored file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the formatting (fonts/fills/borders) of the last existing column
# (J, header + body) for the three new columns K:M, so the new cells look
# consistent with the rest of the table.
$ws.Range("J2:J9").Copy()
$ws.Range("K2:M9").PasteSpecial(-4122)  # xlPasteFormats

# New header captions (shared strings get appended automatically)
$ws.Range("K2").Value = "One tweet alerts"
$ws.Range("L2").Value = "Topics ignoring 1 tweet alerts "
$ws.Range("M2").Value = "Regions ignoring 1 tweet alerts "

# Wrap the (now longer) header captions and grow the header row to fit
$ws.Range("B2:M2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 23.85

# Widen the two new columns to comfortably fit their wrapped captions
$ws.Columns.Item(11).ColumnWidth = 17.09
$ws.Columns.Item(12).ColumnWidth = 20.75

# Slightly wider sheet-wide default column width
$ws.StandardWidth = 11.55078125

# Restore the authored cursor position
[void]$ws.Range("D16").Select()
